# Thurs 17th Dec - import + PEP8 compliance
# Add two new additional-code-type rows to the "Updated" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "Prohibition / restriction /  surveillance"

$ws.Range("A7").Value = 9
$ws.Range("B7").Value = "Export refunds"

# Match the author's final selection (cell A7) recorded in the saved file.
$ws.Range("A7").Select() | Out-Null
